# Update the "Google" worksheet, row 3 (D3:F3) with new filter arguments.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

$ws.Range("D3").Value = "fail"
$ws.Range("E3").Value = "Zsinj | Wookieepedia | FANDOM powered by Wikia"
$ws.Range("F3").Value = "https://starwars.fandom.com/wiki/Zsinj"
